$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "x"
$ws.Range("A2").Value = 10

[void]$ws.Range("A2").Select()
